$d = $word.ActiveDocument
$rng = $d.Content

$rng.Find.Execute("2023-04-13 Thursday", $true, $false, $false, $false, $false, $true, 1, $false, "2023-04-14 Friday", 2) | Out-Null
$rng.Find.Execute("10×97=", $true, $false, $false, $false, $false, $true, 1, $false, "58×77=", 2) | Out-Null
$rng.Find.Execute("78×20=", $true, $false, $false, $false, $false, $true, 1, $false, "18×70=", 2) | Out-Null
$rng.Find.Execute("64×30=", $true, $false, $false, $false, $false, $true, 1, $false, "50×98=", 2) | Out-Null
$rng.Find.Execute("60×83=", $true, $false, $false, $false, $false, $true, 1, $false, "40×99=", 2) | Out-Null
$rng.Find.Execute("84×64=", $true, $false, $false, $false, $false, $true, 1, $false, "92×84=", 2) | Out-Null
$rng.Find.Execute("90×99=", $true, $false, $false, $false, $false, $true, 1, $false, "94×41=", 2) | Out-Null
$rng.Find.Execute("87×91=", $true, $false, $false, $false, $false, $true, 1, $false, "63×72=", 2) | Out-Null
$rng.Find.Execute("40×43=", $true, $false, $false, $false, $false, $true, 1, $false, "85×88=", 2) | Out-Null
$rng.Find.Execute("35×17=", $true, $false, $false, $false, $false, $true, 1, $false, "39×20=", 2) | Out-Null
$rng.Find.Execute("92×76=", $true, $false, $false, $false, $false, $true, 1, $false, "59×31=", 2) | Out-Null
$rng.Find.Execute("91×84=", $true, $false, $false, $false, $false, $true, 1, $false, "52×91=", 2) | Out-Null
$rng.Find.Execute("28×18=", $true, $false, $false, $false, $false, $true, 1, $false, "53×95=", 2) | Out-Null
$rng.Find.Execute("29×19=", $true, $false, $false, $false, $false, $true, 1, $false, "34×77=", 2) | Out-Null
$rng.Find.Execute("15×98=", $true, $false, $false, $false, $false, $true, 1, $false, "38×94=", 2) | Out-Null
$rng.Find.Execute("84×62=", $true, $false, $false, $false, $false, $true, 1, $false, "44×91=", 2) | Out-Null
$rng.Find.Execute("69×17=", $true, $false, $false, $false, $false, $true, 1, $false, "87×38=", 2) | Out-Null
$rng.Find.Execute("23×93=", $true, $false, $false, $false, $false, $true, 1, $false, "37×89=", 2) | Out-Null
$rng.Find.Execute("99×57=", $true, $false, $false, $false, $false, $true, 1, $false, "57×46=", 2) | Out-Null
$rng.Find.Execute("14×81=", $true, $false, $false, $false, $false, $true, 1, $false, "84×11=", 2) | Out-Null
$rng.Find.Execute("51×57=", $true, $false, $false, $false, $false, $true, 1, $false, "37×67=", 2) | Out-Null
$rng.Find.Execute("90×74=", $true, $false, $false, $false, $false, $true, 1, $false, "43×23=", 2) | Out-Null
$rng.Find.Execute("58×97=", $true, $false, $false, $false, $false, $true, 1, $false, "86×34=", 2) | Out-Null
$rng.Find.Execute("15×51=", $true, $false, $false, $false, $false, $true, 1, $false, "94×20=", 2) | Out-Null
$rng.Find.Execute("69×81=", $true, $false, $false, $false, $false, $true, 1, $false, "34×54=", 2) | Out-Null
$rng.Find.Execute("57×67=", $true, $false, $false, $false, $false, $true, 1, $false, "62×30=", 2) | Out-Null
$rng.Find.Execute("95×52=", $true, $false, $false, $false, $false, $true, 1, $false, "69×24=", 2) | Out-Null
$rng.Find.Execute("24×51=", $true, $false, $false, $false, $false, $true, 1, $false, "100×10=", 2) | Out-Null
$rng.Find.Execute("61×50=", $true, $false, $false, $false, $false, $true, 1, $false, "91×87=", 2) | Out-Null
$rng.Find.Execute("80×63=", $true, $false, $false, $false, $false, $true, 1, $false, "42×42=", 2) | Out-Null
$rng.Find.Execute("53×34=", $true, $false, $false, $false, $false, $true, 1, $false, "46×73=", 2) | Out-Null
$rng.Find.Execute("47×22=", $true, $false, $false, $false, $false, $true, 1, $false, "48×11=", 2) | Out-Null
$rng.Find.Execute("30×36=", $true, $false, $false, $false, $false, $true, 1, $false, "17×100=", 2) | Out-Null
$rng.Find.Execute("53×97=", $true, $false, $false, $false, $false, $true, 1, $false, "70×33=", 2) | Out-Null
$rng.Find.Execute("10×16=", $true, $false, $false, $false, $false, $true, 1, $false, "74×62=", 2) | Out-Null
$rng.Find.Execute("86×10=", $true, $false, $false, $false, $false, $true, 1, $false, "23×85=", 2) | Out-Null
$rng.Find.Execute("99×59=", $true, $false, $false, $false, $false, $true, 1, $false, "93×73=", 2) | Out-Null
$rng.Find.Execute("65×70=", $true, $false, $false, $false, $false, $true, 1, $false, "79×47=", 2) | Out-Null
$rng.Find.Execute("25×96=", $true, $false, $false, $false, $false, $true, 1, $false, "33×28=", 2) | Out-Null
$rng.Find.Execute("16×67=", $true, $false, $false, $false, $false, $true, 1, $false, "41×64=", 2) | Out-Null
$rng.Find.Execute("14×79=", $true, $false, $false, $false, $false, $true, 1, $false, "80×38=", 2) | Out-Null
$rng.Find.Execute("32×87=", $true, $false, $false, $false, $false, $true, 1, $false, "58×42=", 2) | Out-Null
$rng.Find.Execute("14×51=", $true, $false, $false, $false, $false, $true, 1, $false, "71×37=", 2) | Out-Null
$rng.Find.Execute("92×90=", $true, $false, $false, $false, $false, $true, 1, $false, "89×20=", 2) | Out-Null
$rng.Find.Execute("10×22=", $true, $false, $false, $false, $false, $true, 1, $false, "15×53=", 2) | Out-Null
$rng.Find.Execute("14×35=", $true, $false, $false, $false, $false, $true, 1, $false, "31×14=", 2) | Out-Null
$rng.Find.Execute("65×56=", $true, $false, $false, $false, $false, $true, 1, $false, "99×45=", 2) | Out-Null
$rng.Find.Execute("38×65=", $true, $false, $false, $false, $false, $true, 1, $false, "59×52=", 2) | Out-Null
$rng.Find.Execute("75×81=", $true, $false, $false, $false, $false, $true, 1, $false, "15×17=", 2) | Out-Null
$rng.Find.Execute("29×57=", $true, $false, $false, $false, $false, $true, 1, $false, "43×83=", 2) | Out-Null
$rng.Find.Execute("99×97=", $true, $false, $false, $false, $false, $true, 1, $false, "41×18=", 2) | Out-Null
$rng.Find.Execute("53×22=", $true, $false, $false, $false, $false, $true, 1, $false, "85×35=", 2) | Out-Null
$rng.Find.Execute("41×17=", $true, $false, $false, $false, $false, $true, 1, $false, "77×29=", 2) | Out-Null
$rng.Find.Execute("95×88=", $true, $false, $false, $false, $false, $true, 1, $false, "54×23=", 2) | Out-Null
$rng.Find.Execute("10×20=", $true, $false, $false, $false, $false, $true, 1, $false, "45×34=", 2) | Out-Null
$rng.Find.Execute("39×54=", $true, $false, $false, $false, $false, $true, 1, $false, "47×67=", 2) | Out-Null
$rng.Find.Execute("77×15=", $true, $false, $false, $false, $false, $true, 1, $false, "39×66=", 2) | Out-Null
$rng.Find.Execute("95×63=", $true, $false, $false, $false, $false, $true, 1, $false, "55×64=", 2) | Out-Null
$rng.Find.Execute("61×29=", $true, $false, $false, $false, $false, $true, 1, $false, "95×74=", 2) | Out-Null
$rng.Find.Execute("14×34=", $true, $false, $false, $false, $false, $true, 1, $false, "90×20=", 2) | Out-Null
$rng.Find.Execute("60×96=", $true, $false, $false, $false, $false, $true, 1, $false, "83×31=", 2) | Out-Null
$rng.Find.Execute("55×58=", $true, $false, $false, $false, $false, $true, 1, $false, "72×29=", 2) | Out-Null
$rng.Find.Execute("39×32=", $true, $false, $false, $false, $false, $true, 1, $false, "85×94=", 2) | Out-Null
$rng.Find.Execute("35×31=", $true, $false, $false, $false, $false, $true, 1, $false, "86×32=", 2) | Out-Null
$rng.Find.Execute("90×66=", $true, $false, $false, $false, $false, $true, 1, $false, "36×57=", 2) | Out-Null
$rng.Find.Execute("31×20=", $true, $false, $false, $false, $false, $true, 1, $false, "65×47=", 2) | Out-Null
$rng.Find.Execute("23×44=", $true, $false, $false, $false, $false, $true, 1, $false, "23×58=", 2) | Out-Null
$rng.Find.Execute("23×69=", $true, $false, $false, $false, $false, $true, 1, $false, "17×71=", 2) | Out-Null
$rng.Find.Execute("43×85=", $true, $false, $false, $false, $false, $true, 1, $false, "29×12=", 2) | Out-Null
$rng.Find.Execute("22×88=", $true, $false, $false, $false, $false, $true, 1, $false, "30×72=", 2) | Out-Null
$rng.Find.Execute("67×39=", $true, $false, $false, $false, $false, $true, 1, $false, "59×96=", 2) | Out-Null
$rng.Find.Execute("74×84=", $true, $false, $false, $false, $false, $true, 1, $false, "88×69=", 2) | Out-Null
$rng.Find.Execute("86×35=", $true, $false, $false, $false, $false, $true, 1, $false, "27×80=", 2) | Out-Null
$rng.Find.Execute("99×93=", $true, $false, $false, $false, $false, $true, 1, $false, "49×68=", 2) | Out-Null
$rng.Find.Execute("52×61=", $true, $false, $false, $false, $false, $true, 1, $false, "25×22=", 2) | Out-Null
$rng.Find.Execute("27×18=", $true, $false, $false, $false, $false, $true, 1, $false, "14×51=", 2) | Out-Null
$rng.Find.Execute("65×98=", $true, $false, $false, $false, $false, $true, 1, $false, "68×88=", 2) | Out-Null
$rng.Find.Execute("85×58=", $true, $false, $false, $false, $false, $true, 1, $false, "27×76=", 2) | Out-Null
$rng.Find.Execute("92×54=", $true, $false, $false, $false, $false, $true, 1, $false, "22×90=", 2) | Out-Null
$rng.Find.Execute("99×19=", $true, $false, $false, $false, $false, $true, 1, $false, "61×55=", 2) | Out-Null
$rng.Find.Execute("65×83=", $true, $false, $false, $false, $false, $true, 1, $false, "59×41=", 2) | Out-Null
$rng.Find.Execute("88×16=", $true, $false, $false, $false, $false, $true, 1, $false, "75×22=", 2) | Out-Null
$rng.Find.Execute("49×65=", $true, $false, $false, $false, $false, $true, 1, $false, "59×14=", 2) | Out-Null
$rng.Find.Execute("93×54=", $true, $false, $false, $false, $false, $true, 1, $false, "63×16=", 2) | Out-Null
$rng.Find.Execute("23×56=", $true, $false, $false, $false, $false, $true, 1, $false, "11×86=", 2) | Out-Null
$rng.Find.Execute("35×83=", $true, $false, $false, $false, $false, $true, 1, $false, "56×65=", 2) | Out-Null
$rng.Find.Execute("47×88=", $true, $false, $false, $false, $false, $true, 1, $false, "54×21=", 2) | Out-Null
$rng.Find.Execute("57×57=", $true, $false, $false, $false, $false, $true, 1, $false, "48×19=", 2) | Out-Null
$rng.Find.Execute("58×40=", $true, $false, $false, $false, $false, $true, 1, $false, "70×79=", 2) | Out-Null
$rng.Find.Execute("87×77=", $true, $false, $false, $false, $false, $true, 1, $false, "34×72=", 2) | Out-Null
$rng.Find.Execute("21×58=", $true, $false, $false, $false, $false, $true, 1, $false, "47×25=", 2) | Out-Null
$rng.Find.Execute("42×66=", $true, $false, $false, $false, $false, $true, 1, $false, "32×21=", 2) | Out-Null
$rng.Find.Execute("39×18=", $true, $false, $false, $false, $false, $true, 1, $false, "54×24=", 2) | Out-Null
$rng.Find.Execute("73×93=", $true, $false, $false, $false, $false, $true, 1, $false, "56×24=", 2) | Out-Null
$rng.Find.Execute("98×63=", $true, $false, $false, $false, $false, $true, 1, $false, "95×71=", 2) | Out-Null
$rng.Find.Execute("28×41=", $true, $false, $false, $false, $false, $true, 1, $false, "12×11=", 2) | Out-Null
$rng.Find.Execute("22×23=", $true, $false, $false, $false, $false, $true, 1, $false, "16×63=", 2) | Out-Null
$rng.Find.Execute("58×53=", $true, $false, $false, $false, $false, $true, 1, $false, "71×59=", 2) | Out-Null
$rng.Find.Execute("28×65=", $true, $false, $false, $false, $false, $true, 1, $false, "65×71=", 2) | Out-Null
$rng.Find.Execute("64×28=", $true, $false, $false, $false, $false, $true, 1, $false, "81×36=", 2) | Out-Null

Write-Output "replaced 100 items"
